$wb = $excel.ActiveWorkbook

# Update the Users list: replace "James Craven" with "Thomas Bailey"
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Thomas Bailey"

# Make "Users" the active/selected sheet (was "Contact"), with A2 selected
$usersSheet.Activate()
$usersSheet.Range("A2").Select()
